$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: fill in the 5th homework score (column G) for this student
$ws.Range("G9").Value = 5

# Row 25: student completed the first two homeworks (columns C and D),
# so their scores become 5 and the "not done" (green) fill is cleared
# to match the normal border style used by already-graded cells.
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5

$ws.Range("G25").Copy()
$ws.Range("C25:D25").PasteSpecial(-4122)
